$wb = $excel.ActiveWorkbook

# Locate the sheet after which the new sheet "tc048" must be inserted (tc029),
# and the sheet that currently follows it (tc007), so we can place the
# new sheet directly between them.
$wsAfter = $wb.Sheets.Item("tc007")

# Add a new worksheet before "tc007" -> this puts it right after "tc029".
$newSheet = $wb.Worksheets.Add($wsAfter)
$newSheet.Name = "tc048"

# Populate header row and data row.
$newSheet.Range("A1").Value = "ModuleName"
$newSheet.Range("B1").Value = "ReqId"
$newSheet.Range("C1").Value = "user"

$newSheet.Range("A2").Value = "Epic Mohit"
$newSheet.Range("B2").Value = "RQ-463"
$newSheet.Range("C2").Value = "Mohit Aman"

# Make the new sheet the active / selected tab.
$newSheet.Activate()
$newSheet.Range("C7").Select()
